$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'34.414.04"
$ws.Range("E2").Value = '  -0.03%  '

# Row 3
$ws.Range("D3").Value = "'1.802.90"
$ws.Range("E3").Value = '  -0.01%  '

# Row 4
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").Value = "'225.24"
$ws.Range("E5").Value = '  -0.98%  '

# Row 6
$ws.Range("D6").Value = "'0.602"
$ws.Range("E6").Value = '  +3.88%  '

# Row 7
$ws.Range("E7").Value = '  -0.14%  '

# Row 8
$ws.Range("D8").Value = "'36.19"
$ws.Range("E8").Value = '  +3.43%  '

# Row 9
$ws.Range("E9").Value = '  -2.22%  '

# Row 10
$ws.Range("D10").Value = "'0.0677"
$ws.Range("E10").Value = '  -1.97%  '

# Row 11
$ws.Range("D11").Value = "'0.0966"
$ws.Range("E11").Value = '  +1.56%  '

# Row 12
$ws.Range("D12").Value = "'2.062.76"
$ws.Range("E12").Value = '  -0.02%  '

# Row 13
$ws.Range("D13").Value = "'11.25"
$ws.Range("E13").Value = '  +0.23%  '

# Row 14
$ws.Range("D14").Value = "'1.801.42"
$ws.Range("E14").Value = '  -0.01%  '

# Row 15
$ws.Range("D15").Value = "'0.628"
$ws.Range("E15").Value = '  -1.92%  '

# Row 16
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = "'34.385.91"
$ws.Range("E16").Value = '  -0.03%  '

# Row 17
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").Value = "'4.44"
$ws.Range("E17").Value = '  +2.40%  '

# Row 18
$ws.Range("D18").Value = "'68.56"
$ws.Range("E18").Value = '  -0.74%  '

# Row 19
$ws.Range("D19").Value = "'242.53"
$ws.Range("E19").Value = '  -0.92%  '

# Row 21
$ws.Range("D21").Value = "'11.24"
$ws.Range("E21").Value = '  -2.66%  '

# Row 22
$ws.Range("E22").Value = '  -0.14%  '

# Row 23
$ws.Range("E23").Value = '  -1.54%  '

# Row 24
$ws.Range("E24").Value = '  +4.87%  '

# Row 25
$ws.Range("D25").Value = "'170.58"
$ws.Range("E25").Value = '  +0.20%  '

# Row 26
$ws.Range("D26").Value = "'7.88"
$ws.Range("E26").Value = '  +4.22%  '

# Row 27
$ws.Range("D27").Value = "'17.44"
$ws.Range("E27").Value = '  +4.25%  '

# Row 28
$ws.Range("E28").Value = '  +2.37%  '

# Row 29
$ws.Range("E29").Value = '  -0.15%  '

# Row 30
$ws.Range("B30").Value = 'InternetComputer(DFINITY)'
$ws.Range("C30").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D30").Value = "'3.92"
$ws.Range("E30").Value = '  -1.70%  '

# Row 31
$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").Value = "'3.79"
$ws.Range("E31").Value = '  -0.30%  '

# Row 32
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").Value = "'1.23"
$ws.Range("E32").Value = '  -1.25%  '

# Row 33
$ws.Range("D33").Value = "'0.0515"
$ws.Range("E33").Value = '  -2.38%  '

# Row 34
$ws.Range("E34").Value = '  -2.99%  '

# Row 35
$ws.Range("D35").Value = "'1.363.16"
$ws.Range("E35").Value = '  -2.61%  '

# Row 36
$ws.Range("D36").Value = "'0.650"
$ws.Range("E36").Value = '  -4.51%  '

# Row 37
$ws.Range("E37").Value = '  -0.37%  '

# Row 38
$ws.Range("D38").Value = "'2.36"
$ws.Range("E38").Value = '  -6.92%  '

# Row 39
$ws.Range("E39").Value = '  -1.81%  '

# Row 40
$ws.Range("E40").Value = '  +1.48%  '

# Row 41
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = '  -1.79%  '

# Row 42
$ws.Range("D42").Value = "'80.84"
$ws.Range("E42").Value = '  -2.73%  '

# Row 43
$ws.Range("D43").Value = "'0.937"
$ws.Range("E43").Value = '  -0.94%  '

# Row 44
$ws.Range("E44").Value = '  +5.15%  '

# Row 45
$ws.Range("D45").Value = "'13.27"
$ws.Range("E45").Value = '  -2.77%  '

# Row 46
$ws.Range("E46").Value = '  -2.89%  '

# Row 47
$ws.Range("D47").Value = "'1.964.62"
$ws.Range("E47").Value = '  +0.05%  '

# Row 48
$ws.Range("D48").Value = "'5.78"
$ws.Range("E48").Value = '  -3.46%  '

# Row 49
$ws.Range("E49").Value = '  -0.16%  '

# Row 50
$ws.Range("D50").Value = "'102.16"
$ws.Range("E50").Value = '  -2.09%  '

# Row 51
$dVal = "'{0}{1}{2}" -f '0.0', [char]8326, '0124'
$ws.Range("D51").Value = $dVal
$ws.Range("E51").Value = '  -3.66%  '
